# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) pair recorded on rows 16
# and 17 were swapped: the period that used to be on row 16 (1811 / 34200)
# moved to row 17, and the period that used to be on row 17 (1810 / 38000)
# moved to row 16.
#
# NOTE: use Value2 (not Value) to read/write cell contents - in this
# COM-interop runtime, Range.Value resolves to the raw Variant property
# accessor instead of invoking it, so Value2 is used for reliable get/set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$e16 = $ws.Range("E16").Value2
$e17 = $ws.Range("E17").Value2
$f16 = $ws.Range("F16").Value2
$f17 = $ws.Range("F17").Value2

$ws.Range("E16").Value2 = $e17
$ws.Range("E17").Value2 = $e16
$ws.Range("F16").Value2 = $f17
$ws.Range("F17").Value2 = $f16
